$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from H1 onto the two new header cells I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Header cells (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows (plain numeric values, no special style)
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
